$wb = $excel.ActiveWorkbook

# --- Capture references to the existing sheets before inserting the new one ---
$zongji  = $wb.Worksheets.Item(1)     # "总计"
$oldQ3   = $wb.Worksheets.Item(2)     # "2022-Q3" (index shifts later, object ref stays valid)
$oldQ2   = $wb.Worksheets.Item(3)     # "2022-Q2"
$oldQ1   = $wb.Worksheets.Item(4)     # "2022-Q1"
$old21Q4 = $wb.Worksheets.Item(5)     # "2021-Q4"

# --- Insert the new "2022-Q4" sheet right after "总计" ---
$q4 = $wb.Worksheets.Add($null, $zongji)
$q4.Name = "2022-Q4"

# ============================================================
# 1. Update the "总计" (summary) sheet: insert a new data row
#    for 2022-Q4 at the top, shifting the rest down by one row.
# ============================================================

# Shift existing rows 5->6, 4->5, 3->4, 2->3 (copy full formatting + values)
for ($r = 5; $r -ge 2; $r--) {
    $dstRow = $r + 1
    $zongji.Range("A$r" + ":D$r").Copy($zongji.Range("A$dstRow" + ":D$dstRow"))
    $zongji.Range("A$dstRow").Value = $r - 1
}

# Write the new 2022-Q4 summary row into row 2 (keeps row 2's existing style)
$zongji.Range("A2").Value = 0
$zongji.Range("B2").Value = "2022-Q4"
$zongji.Range("C2").Value = 27
$zongji.Range("D2").Value = 2.52

# ============================================================
# 2. Populate the new "2022-Q4" fund-holdings sheet.
#    Template formatting (header + per-row styles) is copied
#    from the existing "2022-Q2" sheet, which already has the
#    same header / column layout used by every quarter sheet.
# ============================================================

$fundData = @(
    @(0, "000478", "建信中证500指数增强A", "48.60", "83.82", "0.95", "0.4617", 9),
    @(1, "001556", "天弘中证500指数增强A", "25.50", "94.27", "1.74", "0.4437", 4),
    @(2, "007994", "华夏中证500指数增强A", "22.71", "93.73", "1.10", "0.2498", 6),
    @(3, "013233", "华夏中证500指数智选增强A", "21.33", "93.95", "1.14", "0.2432", 4),
    @(4, "001557", "天弘中证500指数增强C", "13.20", "94.27", "1.74", "0.2297", 4),
    @(5, "005994", "国投瑞银中证500指数量化增强A", "13.06", "89.93", "1.58", "0.2063", 2),
    @(6, "013641", "博道成长智航股票A", "9.90", "91.40", "1.02", "0.1010", 6),
    @(7, "159610", "景顺长城中证500增强策略ETF", "6.09", "98.72", "1.18", "0.0719", 5),
    @(8, "013642", "博道成长智航股票C", "6.95", "91.40", "1.02", "0.0709", 6),
    @(9, "007995", "华夏中证500指数增强C", "6.02", "93.73", "1.10", "0.0662", 6),
    @(10, "007089", "国投瑞银中证500指数量化增强C", "3.73", "89.93", "1.58", "0.0589", 2),
    @(11, "005062", "博时中证500指数增强A", "3.43", "93.43", "1.61", "0.0552", 9),
    @(12, "013234", "华夏中证500指数智选增强C", "4.20", "93.95", "1.14", "0.0479", 4),
    @(13, "162216", "泰达宏利中证500指数增强（LOF）", "2.78", "93.77", "1.22", "0.0339", 7),
    @(14, "005633", "建信中证500指数增强C", "3.38", "83.82", "0.95", "0.0321", 9),
    @(15, "159804", "国寿安保国证创业板中盘精选88ETF", "1.15", "99.00", "1.83", "0.0210", 7),
    @(16, "006441", "中信建投中证500指数增强C", "2.27", "93.50", "0.87", "0.0197", 4),
    @(17, "006440", "中信建投中证500指数增强A", "2.19", "93.50", "0.87", "0.0191", 4),
    @(18, "005795", "博时中证500指数增强C", "1.12", "93.43", "1.61", "0.0180", 9),
    @(19, "014344", "鹏华中证500指数增强A", "0.99", "92.67", "1.81", "0.0179", 3),
    @(20, "000270", "建信灵活配置混合", "1.53", "91.22", "1.03", "0.0158", 2),
    @(21, "014345", "鹏华中证500指数增强C", "0.86", "92.67", "1.81", "0.0156", 3),
    @(22, "012498", "汇添富中证500基本面增强指数A", "0.77", "92.29", "1.94", "0.0149", 4),
    @(23, "012499", "汇添富中证500基本面增强指数C", "0.23", "92.29", "1.94", "0.0045", 4),
    @(24, "005260", "银华稳健增利灵活配置混合A", "0.33", "91.18", "0.65", "0.0021", 10),
    @(25, "005261", "银华稳健增利灵活配置混合C", "0.21", "91.18", "0.65", "0.0014", 10),
    @(26, "015245", "南华丰汇混合", "0.11", "84.24", "1.07", "0.0012", 8)
)

$rowCount = $fundData.Count

# Copy header row formatting + text
$oldQ2.Range("A1:H1").Copy($q4.Range("A1:H1"))

# Copy one data-row's formatting down for every fund row we need
for ($i = 0; $i -lt $rowCount; $i++) {
    $r = 2 + $i
    $oldQ2.Range("A2:H2").Copy($q4.Range("A$r" + ":H$r"))
}

# Write column A (index) and column H (rank) as real numbers -- these keep
# the per-row style/format that was just copied above.
for ($i = 0; $i -lt $rowCount; $i++) {
    $r = 2 + $i
    $row = $fundData[$i]
    $q4.Range("A$r").Value = $row[0]
    $q4.Range("H$r").Value = $row[7]
}

# Columns B-G must stay plain text (fund codes with leading zeros, and
# numeric-looking strings like "48.60" that must not collapse to 48.6).
# Stage the text in a scratch area far to the right, force it to Text
# format there, copy it, and paste-special (values only) onto the real
# destination so the destination's copied style/format is left untouched.
$scratchCol0 = 50   # column AX, well clear of the real data (A-H)
for ($i = 0; $i -lt $rowCount; $i++) {
    $r = 2 + $i
    $row = $fundData[$i]
    for ($c = 0; $c -lt 6; $c++) {
        $cell = $q4.Cells.Item($r, $scratchCol0 + $c)
        $cell.NumberFormat = "@"
        $cell.Value = $row[$c + 1]
    }
}
$scratchRange = $q4.Range($q4.Cells.Item(2, $scratchCol0), $q4.Cells.Item(1 + $rowCount, $scratchCol0 + 5))
$scratchRange.Copy()
$destRange = $q4.Range($q4.Cells.Item(2, 2), $q4.Cells.Item(1 + $rowCount, 7))
$destRange.PasteSpecial(-4163)
$scratchRange.Clear()

Write-Host "edit complete"
